# "Generate Report for Archive"
# Update localization status from "Ready for handoff" to "In Translation"
# across the Overview summary sheet and each per-locale detail sheet, then
# let the column widths re-flow to fit the (now shorter) status text.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: per-locale status columns (E = zh-cn, F = de-de) for each
# tracked file row.
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus

# Per-locale detail sheets: Status column (C) for each tracked file row.
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("C3").Value = $newStatus

$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("C3").Value = $newStatus

# The status columns auto-size to the new, shorter text.
$ws_overview.Columns.Item(5).ColumnWidth = 12.5
$ws_overview.Columns.Item(6).ColumnWidth = 12.5
$ws_zhcn.Columns.Item(3).ColumnWidth = 12.5
$ws_dede.Columns.Item(3).ColumnWidth = 12.5
